# Generate Report for Handoff
# Updates the localization-status workbook after a handoff report
# regeneration: six rows (a7bcba37..md through fb3553fe..md, skipping the
# already-localized e1feb996..md row) now carry a refreshed handoff
# timestamp and a "ht" priority flag on the per-language sheets.

$wb = $excel.ActiveWorkbook

$rows = @(8, 9, 10, 11, 12, 14)

# Overview sheet: "Latest HO Xliff Generate Date" column G
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-09-02 10:26:13"
}

# zh-cn sheet: "Latest Handoff Datetime" column H, "Priority" column E
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-09-02 10:26:00"
}

# de-de sheet: "Latest Handoff Datetime" column H, "Priority" column E
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-09-02 10:26:13"
}
